$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename link_* headers (row 1), columns U..AW ---
# executionPlan -> executionDetails, jobExecutionDetails -> executionLinks,
# testProjectNavigation -> projectLinks, testRunNavigation -> testRunLinks
$ws.Range("U1").Value = "link_executionDetails_internalRoleLinkName"
$ws.Range("V1").Value = "link_executionDetails_internalRoleLinkName_1"
$ws.Range("W1").Value = "link_executionDetails_job_id"
$ws.Range("X1").Value = "link_executionDetails_job_id_1"
$ws.Range("Y1").Value = "link_executionDetails_plan_id"
$ws.Range("Z1").Value = "link_executionDetails_plan_id_1"
$ws.Range("AA1").Value = "link_executionDetails_project_id"
$ws.Range("AB1").Value = "link_executionDetails_project_id_1"
$ws.Range("AC1").Value = "link_executionDetails_team_id"
$ws.Range("AD1").Value = "link_executionDetails_team_id_1"
$ws.Range("AE1").Value = "link_executionLinks_executions_id"
$ws.Range("AF1").Value = "link_executionLinks_executions_id_1"
$ws.Range("AG1").Value = "link_executionLinks_internalRoleLinkName"
$ws.Range("AH1").Value = "link_executionLinks_internalRoleLinkName_1"
$ws.Range("AI1").Value = "link_executionLinks_project_id"
$ws.Range("AJ1").Value = "link_executionLinks_project_id_1"
$ws.Range("AK1").Value = "link_executionLinks_team_id"
$ws.Range("AL1").Value = "link_executionLinks_team_id_1"
$ws.Range("AM1").Value = "link_projectLinks_internalRoleLinkName"
$ws.Range("AN1").Value = "link_projectLinks_project_id"
$ws.Range("AO1").Value = "link_projectLinks_team_id"
$ws.Range("AP1").Value = "link_projectLinks_test_project_id"
$ws.Range("AQ1").Value = "link_projectLinks_trNthChild"
$ws.Range("AR1").Value = "link_testRunLinks_plan_id"
$ws.Range("AS1").Value = "link_testRunLinks_plan_id_1"
$ws.Range("AT1").Value = "link_testRunLinks_project_id"
$ws.Range("AU1").Value = "link_testRunLinks_project_id_1"
$ws.Range("AV1").Value = "link_testRunLinks_team_id"
$ws.Range("AW1").Value = "link_testRunLinks_team_id_1"

# --- Update row 2 data values to match the re-shuffled link columns ---
# Cells hold numeric-looking text; force text storage (matches the
# original inlineStr/string cell type) via Text format, then clear the
# format again so the cell style index is unaffected.
function Set-TextValue($addr, $val) {
  $r = $ws.Range($addr)
  $r.NumberFormat = "@"
  $r.Value = $val
  $r.ClearFormats()
}

Set-TextValue "U2" "8"
Set-TextValue "V2" "10"
Set-TextValue "W2" "8"
Set-TextValue "X2" "10"
Set-TextValue "Y2" "837097"
Set-TextValue "Z2" "837132"
Set-TextValue "AA2" "1588984"
Set-TextValue "AB2" "1588984"
Set-TextValue "AC2" "1570311"
Set-TextValue "AD2" "1570311"
Set-TextValue "AE2" "10"
Set-TextValue "AF2" "12"
Set-TextValue "AG2" "10"
Set-TextValue "AH2" "12"

# --- Column widths follow the renamed headers (columns U..AW, 21..49) ---
# Width (in characters) = len(header text) + 2, matching this sheet's
# existing convention throughout (see columns A..AX). Assigning
# ColumnWidth directly re-applies a +5/6 character padding on save, so
# subtract it up front to land on the clean integer width the workbook
# already uses everywhere else.
$newWidths = @{
  21 = 44; 22 = 46; 23 = 30; 24 = 32; 25 = 31; 26 = 33; 27 = 34; 28 = 36;
  29 = 31; 30 = 33; 31 = 35; 32 = 37; 33 = 42; 34 = 44; 35 = 32; 36 = 34;
  37 = 29; 38 = 31; 39 = 40; 40 = 30; 41 = 27; 42 = 35; 43 = 30; 44 = 27;
  45 = 29; 46 = 30; 47 = 32; 48 = 27; 49 = 29
}
foreach ($col in $newWidths.Keys) {
  $ws.Columns.Item($col).ColumnWidth = ($newWidths[$col] - 5.0/6.0)
}
